# Edit script for Ex27-Pi2GoSimulator-Objects.docx
#
# Changes applied (per the commit "Initio Exercises 10. Adding aims to
# all exercise sheets."):
#   1. A new "AIM: ..." paragraph (plus a following blank paragraph) is
#      inserted immediately before the "Exercise 1:" paragraph. The
#      "_GoBack" bookmark that used to sit after "Hint" now sits inside
#      this new AIM paragraph (mid-word, between "Obj" and "ects").
#   2. The old "_GoBack" bookmark next to "Hint" is removed (it moved,
#      see above).
#   3. The two runs making up " University of Liverpool, 20" + "20" are
#      merged into a single run " University of Liverpool, 2020".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: remove the old "_GoBack" bookmark (next to "Hint") first, so
# there is no name clash with the one we insert in step 2.
# ---------------------------------------------------------------------
$bookmarks = $d.Bookmarks
if ($bookmarks.Exists("_GoBack")) {
    $bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# Step 2: insert the new "AIM:" paragraph + a following blank paragraph
# right before the "Exercise 1:" paragraph, using a literal OOXML
# fragment so the paragraph/run formatting matches exactly (explicit
# Calibri fonts, no theme refs) and the bookmark lands mid-run exactly
# where it used to be relative to the text ("Obj" | "ects").
# ---------------------------------------------------------------------
$exerciseParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Exercise 1:")) {
        $exerciseParaIndex = $i
        break
    }
}
if ($exerciseParaIndex -eq -1) {
    throw "Could not locate the 'Exercise 1:' paragraph"
}

$insertPos = $d.Paragraphs($exerciseParaIndex).Range.Start
$insertRange = $d.Range($insertPos, $insertPos)

$aimXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/></w:rPr><w:t xml:space="preserve">AIM: </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve">This exercise sheet provides additional </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>exercises using Obj</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>ects and Classes.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertRange.InsertXML($aimXml) | Out-Null

# ---------------------------------------------------------------------
# Step 3: merge the " University of Liverpool, 20" / "20" runs into a
# single run with the text " University of Liverpool, 2020".
# ---------------------------------------------------------------------
$univParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.Contains("University of Liverpool")) {
        $univParaIndex = $i
        break
    }
}
if ($univParaIndex -eq -1) {
    throw "Could not locate the 'University of Liverpool' paragraph"
}

$univPara = $d.Paragraphs($univParaIndex)
$splitPos = $univPara.Range.Start + [string](" University of Liverpool, 20").Length
$trailingRun = $d.Range($splitPos, $splitPos + 2)
if ($trailingRun.Text -eq "20") {
    $trailingRun.Delete()
    $gap = $d.Range($splitPos, $splitPos)
    $gap.InsertAfter("20")
}

Write-Output "done"
